# Update frontend load data from backend
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subject codes changed from "cn101" to "cn103" for every student row.
$ws.Range("C2:C5").Value = "cn103"

# The "grade" column (D) is dropped entirely; "midterm" (old column E)
# becomes the new column D. Clear D's old contents, then move E into D.
$ws.Range("D1:D5").ClearContents()
$ws.Range("E1:E5").Cut($ws.Range("D1:D5"))

# The now-empty header cell left behind at E1 is removed outright (no
# leftover empty cell), while F1:H1 stay untouched.
$ws.Range("E1").Clear()

# Selection moves to the new "midterm" column.
$ws.Range("D1:D5").Select()
